$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# Extend row 1 (header) and row 2 (data) formatting into the new columns
# H:N by copying the existing header/body cell formats (bold+border for
# the header row, plain for the data row) before writing values into them.
$ws.Range("B1:G1").Copy() | Out-Null
$ws.Range("H1:N1").PasteSpecial(-4122) | Out-Null
$ws.Range("B2:G2").Copy() | Out-Null
$ws.Range("H2:N2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Header row: rename/extend columns to the common property schema used by
# the other sheets (name, capacity, owner, register_date, register_reason,
# acquire_value, property_category, category, date, legislator_name,
# legislator_id, source_file, index).
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Data row 2: fill in the new trailing columns for this land/car property
# record. J2 is written as a formula that evaluates to the literal text
# "2013-12-19" and then flattened back to a value in place; this keeps it
# a plain text cell without Excel coercing the date-looking string into a
# date serial (and without leaving a stray NumberFormat style behind).
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Formula = '="2013-12-19"'
$ws.Range("J2").Copy() | Out-Null
$ws.Range("J2").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("K2").Value = "楊玉欣"
$ws.Range("L2").Value = 1757
$ws.Range("M2").Value = "tmp7d8c1"
$ws.Range("N2").Value = 29
